$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 363.5
